# Apply updated "view count" (column F) values to the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (sheet1.xml) changes
$ws1.Range("F7").Value = 1229
$ws1.Range("F9").Value = 806
$ws1.Range("F13").Value = 355
$ws1.Range("F15").Value = 930
$ws1.Range("F16").Value = 9854
$ws1.Range("F17").Value = 615
$ws1.Range("F33").Value = 97
$ws1.Range("F36").Value = 193
$ws1.Range("F37").Value = 170
$ws1.Range("F38").Value = 39

# Sheet "全部类型" (sheet4.xml) changes
$ws4.Range("F11").Value = 1229
$ws4.Range("F14").Value = 806
$ws4.Range("F17").Value = 355
$ws4.Range("F19").Value = 930
$ws4.Range("F20").Value = 9854
$ws4.Range("F22").Value = 615
$ws4.Range("F40").Value = 97
$ws4.Range("F46").Value = 193
$ws4.Range("F47").Value = 170
